# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg) values
# across the data rows (rows 3-20; row 2 and row 5 are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $D, $M, $N, $O, $P, $S)

    $ws.Cells.Item($Row, 4).Value = $D    # D - Fecha
    $ws.Cells.Item($Row, 13).Value = $M   # M - Volumen
    $ws.Cells.Item($Row, 14).Value = $N   # N - Precio minimo
    $ws.Cells.Item($Row, 15).Value = $O   # O - Precio maximo
    $ws.Cells.Item($Row, 16).Value = $P   # P - Precio promedio ponderado
    $ws.Cells.Item($Row, 19).Value = $S   # S - Precio $/Kg
}

Set-Row 3  44448 30 22000 22000 22000 1100
Set-Row 4  44376 38 20000 20000 20000 1000
Set-Row 6  44298 65 22000 22000 22000 1100
Set-Row 7  44377 25 20000 20000 20000 1000
Set-Row 8  44292 30 25000 25000 25000 1250
Set-Row 9  44406 20 20000 20000 20000 1000
Set-Row 10 44305 20 22000 22000 22000 1100
Set-Row 11 44382 24 20000 20000 20000 1000
Set-Row 12 44291 70 25000 25000 25000 1250
Set-Row 13 44300 45 22000 22000 22000 1100
Set-Row 14 44400 45 20000 20000 20000 1000
Set-Row 15 44403 50 20000 20000 20000 1000
Set-Row 16 44445 45 20000 20000 20000 1000
Set-Row 17 44294 25 25000 25000 25000 1250
Set-Row 18 44301 38 22000 22000 22000 1100
Set-Row 19 44385 36 20000 20000 20000 1000
Set-Row 20 44413 45 20000 20000 20000 1000
